# Generate Report for Handback
# Adds a new handback row for "b31a489c-1daa-4a31-8719-4b4fafd4ebac.md"
# to the Overview, zh-cn and de-de sheets (row 4 on each), wires up the
# corresponding hyperlinks, and grows each sheet's table to include the
# new row.

$wb = $excel.ActiveWorkbook

function Set-HyperlinkStyle($rng) {
    $rng.Font.Underline = 2
    $rng.Font.Color = 15570276
}

function Set-DateStyle($rng) {
    $rng.NumberFormat = "yyyy-mm-dd HH:mm:ss"
}

# ---------------------------------------------------------------------
# Sheet "Overview" - new row 4
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A4").Value = "b31a489c-1daa-4a31-8719-4b4fafd4ebac.md"
$wsOverview.Range("B4").Value = "e2e\b31a489c-1daa-4a31-8719-4b4fafd4ebac.md"
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("E4").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F4").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G4").Value = "2016-08-16 22:44:46"

Set-HyperlinkStyle $wsOverview.Range("B4")
Set-DateStyle $wsOverview.Range("G4")

$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9d3ccab124628d52441746074639db66cfceaaf/e2e/b31a489c-1daa-4a31-8719-4b4fafd4ebac.md", "", "", "e2e\b31a489c-1daa-4a31-8719-4b4fafd4ebac.md")

$wsOverview.ListObjects.Item(1).Resize($wsOverview.Range("A1:G4"))

# ---------------------------------------------------------------------
# Sheet "zh-cn" - new row 4
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A4").Value = "b31a489c-1daa-4a31-8719-4b4fafd4ebac.md"
$wsZhCn.Range("B4").Value = ".md"
$wsZhCn.Range("C4").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("D4").Value = "e2e"
$wsZhCn.Range("E4").Value = "ht"
$wsZhCn.Range("F4").Value = "True"
$wsZhCn.Range("G4").Value = "b31a489c-1daa-4a31-8719-4b4fafd4ebac.9d3ccab124628d52441746074639db66cfceaaf0.zh-cn.xlf"
$wsZhCn.Range("H4").Value = "2016-08-16 22:44:41"
$wsZhCn.Range("I4").Value = "b31a489c-1daa-4a31-8719-4b4fafd4ebac.md"
$wsZhCn.Range("J4").Value = "b31a489c-1daa-4a31-8719-4b4fafd4ebac.9d3ccab124628d52441746074639db66cfceaaf0.zh-cn.xlf"
$wsZhCn.Range("K4").Value = "2016-08-16 22:44:59"
$wsZhCn.Range("L4").Value = ""
$wsZhCn.Range("M4").Value = "True"
$wsZhCn.Range("N4").Value = ""
$wsZhCn.Range("O4").Value = "False"
$wsZhCn.Range("P4").Value = ""

Set-HyperlinkStyle $wsZhCn.Range("A4")
Set-HyperlinkStyle $wsZhCn.Range("I4")
Set-DateStyle $wsZhCn.Range("H4")
Set-DateStyle $wsZhCn.Range("K4")

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9d3ccab124628d52441746074639db66cfceaaf/e2e/b31a489c-1daa-4a31-8719-4b4fafd4ebac.md", "", "", "b31a489c-1daa-4a31-8719-4b4fafd4ebac.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/9d3ccab124628d52441746074639db66cfceaaf0/e2e/b31a489c-1daa-4a31-8719-4b4fafd4ebac.md", "", "", "b31a489c-1daa-4a31-8719-4b4fafd4ebac.md")

$wsZhCn.ListObjects.Item(1).Resize($wsZhCn.Range("A1:P4"))

# ---------------------------------------------------------------------
# Sheet "de-de" - new row 4
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A4").Value = "b31a489c-1daa-4a31-8719-4b4fafd4ebac.md"
$wsDeDe.Range("B4").Value = ".md"
$wsDeDe.Range("C4").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("D4").Value = "e2e"
$wsDeDe.Range("E4").Value = "ht"
$wsDeDe.Range("F4").Value = "True"
$wsDeDe.Range("G4").Value = "b31a489c-1daa-4a31-8719-4b4fafd4ebac.9d3ccab124628d52441746074639db66cfceaaf0.de-de.xlf"
$wsDeDe.Range("H4").Value = "2016-08-16 22:44:46"
$wsDeDe.Range("I4").Value = "b31a489c-1daa-4a31-8719-4b4fafd4ebac.md"
$wsDeDe.Range("J4").Value = "b31a489c-1daa-4a31-8719-4b4fafd4ebac.9d3ccab124628d52441746074639db66cfceaaf0.de-de.xlf"
$wsDeDe.Range("K4").Value = "2016-08-16 22:45:14"
$wsDeDe.Range("L4").Value = ""
$wsDeDe.Range("M4").Value = "True"
$wsDeDe.Range("N4").Value = ""
$wsDeDe.Range("O4").Value = "False"
$wsDeDe.Range("P4").Value = ""

Set-HyperlinkStyle $wsDeDe.Range("A4")
Set-HyperlinkStyle $wsDeDe.Range("I4")
Set-DateStyle $wsDeDe.Range("H4")
Set-DateStyle $wsDeDe.Range("K4")

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9d3ccab124628d52441746074639db66cfceaaf/e2e/b31a489c-1daa-4a31-8719-4b4fafd4ebac.md", "", "", "b31a489c-1daa-4a31-8719-4b4fafd4ebac.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/9d3ccab124628d52441746074639db66cfceaaf0/e2e/b31a489c-1daa-4a31-8719-4b4fafd4ebac.md", "", "", "b31a489c-1daa-4a31-8719-4b4fafd4ebac.md")

$wsDeDe.ListObjects.Item(1).Resize($wsDeDe.Range("A1:P4"))
